$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ValidLogin"

$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "manager"

$ws.Range("B3").Select()
